# BLPG Quarterly Financials - insert latest quarter (period ending 2018-09-29)
# and shift the existing 8 quarters of history one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D; existing D:K data (and formatting) moves to E:L.
$ws.Columns("D:D").Insert()

# 2. The freshly inserted column has no format of its own yet - clone it from the
#    column immediately to its right (which is what used to be column D). Only the
#    three data blocks actually carried a column-D cell before the edit, so restrict
#    the format clone to those rows (the section-header rows must stay untouched).
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3. Populate the new column D with the newest quarter's figures.
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 1000
$ws.Range("D9").Value = 800
$ws.Range("D10").Value = 200
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1400
$ws.Range("D18").Value = -400
$ws.Range("D20").Value = 100
$ws.Range("D21").Value = -200
$ws.Range("D22").Value = 300
$ws.Range("D23").Value = -600
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -600
$ws.Range("D27").Value = -600
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -100
$ws.Range("D33").Value = -600
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -600
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 200
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 100
$ws.Range("D46").Value = 300
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 300
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 600
$ws.Range("D57").Value = 700
$ws.Range("D58").Value = 2300
$ws.Range("D59").Value = 600
$ws.Range("D60").Value = 3600
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 3600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -10000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = -3000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -600
$ws.Range("D83").Value = 0
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -200
$ws.Range("D91").Value = -200
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 400
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 0

# 4. Two historical "Capital Expenditures" cells (now I91/J91, previously H91/I91)
#    were revised from 0 to not-available during this update.
$ws.Range("I91").Value = "NA"
$ws.Range("J91").Value = "NA"
